$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2,8).Value = 118
$ws.Cells.Item(2,9).Value = 118
$ws.Cells.Item(2,11).Value = 118
$ws.Cells.Item(2,13).Value = -5

$ws.Cells.Item(32,8).Value = 2001
$ws.Cells.Item(32,9).Value = 0
$ws.Cells.Item(32,10).Value = 2001
$ws.Cells.Item(32,11).Value = 0
$ws.Cells.Item(32,12).Value = 2001
$ws.Cells.Item(32,13).ClearContents()
$ws.Cells.Item(32,14).Value = -2653

$ws.Cells.Item(33,8).Value = 342.26666
$ws.Cells.Item(33,9).Value = 348.6154
$ws.Cells.Item(33,11).Value = 348.6154
$ws.Cells.Item(33,13).Value = -119.6154

$ws.Cells.Item(38,8).Value = 1314.8462
$ws.Cells.Item(38,9).Value = 86.625
$ws.Cells.Item(38,10).Value = 3280
$ws.Cells.Item(38,11).Value = 259.875
$ws.Cells.Item(38,12).Value = 9840
$ws.Cells.Item(38,13).Value = 112.125
$ws.Cells.Item(38,14).Value = -10584

$ws.Cells.Item(39,8).Value = 277
$ws.Cells.Item(39,9).Value = 242.5
$ws.Cells.Item(39,10).Value = 300
$ws.Cells.Item(39,11).Value = 727.5
$ws.Cells.Item(39,12).Value = 900
$ws.Cells.Item(39,13).Value = -431.5
$ws.Cells.Item(39,14).Value = -1492

$ws.Cells.Item(40,8).Value = 1700
$ws.Cells.Item(40,10).Value = 1700
$ws.Cells.Item(40,12).Value = 1700
$ws.Cells.Item(40,14).Value = -2050

$ws.Cells.Item(70,8).Value = 8825810
$ws.Cells.Item(70,10).Value = 1394.5714
$ws.Cells.Item(70,12).Value = 4183.7142
$ws.Cells.Item(70,14).Value = -4723.7142

$ws.Cells.Item(73,8).Value = 8825810
$ws.Cells.Item(73,10).Value = 1394.5714
$ws.Cells.Item(73,12).Value = 4183.7142
$ws.Cells.Item(73,14).Value = -6055.7142

$ws.Cells.Item(81,8).Value = 36000
$ws.Cells.Item(81,10).Value = 36000
$ws.Cells.Item(81,12).Value = 36000
$ws.Cells.Item(81,14).Value = -37996

$ws.Cells.Item(84,8).Value = 36000
$ws.Cells.Item(84,10).Value = 36000
$ws.Cells.Item(84,12).Value = 108000
$ws.Cells.Item(84,14).Value = -117984

$ws.Cells.Item(100,8).Value = 2186.625
$ws.Cells.Item(100,9).Value = 1248.125
$ws.Cells.Item(100,10).Value = 3125.125
$ws.Cells.Item(100,11).Value = 1248.125
$ws.Cells.Item(100,12).Value = 3125.125
$ws.Cells.Item(100,13).Value = -707.125
$ws.Cells.Item(100,14).Value = -4207.125

$ws.Cells.Item(112,8).Value = 1641.2424
$ws.Cells.Item(112,10).Value = 1679.4193
$ws.Cells.Item(112,12).Value = 5038.257900000001
$ws.Cells.Item(112,14).Value = -7254.257900000001

$ws.Cells.Item(129,8).Value = 1038.95
$ws.Cells.Item(129,9).Value = 797.8
$ws.Cells.Item(129,11).Value = 2393.4
$ws.Cells.Item(129,13).Value = 2606.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(22,8).Value = 20000
$ws.Cells.Item(22,10).Value = 26000
$ws.Cells.Item(22,12).Value = 26000
$ws.Cells.Item(22,14).Value = -26598

$ws.Cells.Item(74,8).Value = 1753.5
$ws.Cells.Item(74,9).Value = 1621
$ws.Cells.Item(74,10).Value = 1939
$ws.Cells.Item(74,11).Value = 1621
$ws.Cells.Item(74,12).Value = 1939
$ws.Cells.Item(74,13).Value = -747
$ws.Cells.Item(74,14).Value = -3687

$ws.Cells.Item(77,8).Value = 1753.5
$ws.Cells.Item(77,9).Value = 1621
$ws.Cells.Item(77,10).Value = 1939
$ws.Cells.Item(77,11).Value = 8105
$ws.Cells.Item(77,12).Value = 9695
$ws.Cells.Item(77,13).Value = -3737
$ws.Cells.Item(77,14).Value = -18431

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(5,8).Value = 10000
$ws.Cells.Item(5,9).Value = 0
$ws.Cells.Item(5,10).Value = 10000
$ws.Cells.Item(5,11).Value = 0
$ws.Cells.Item(5,12).Value = 10000
$ws.Cells.Item(5,13).ClearContents()
$ws.Cells.Item(5,14).Value = -10226

$ws.Cells.Item(22,8).Value = 8663.333000000001
$ws.Cells.Item(22,9).Value = 10316
$ws.Cells.Item(22,10).Value = 400
$ws.Cells.Item(22,11).Value = 10316
$ws.Cells.Item(22,12).Value = 400
$ws.Cells.Item(22,13).Value = -10143
$ws.Cells.Item(22,14).Value = -746

$ws.Cells.Item(99,8).Value = 1777.1904
$ws.Cells.Item(99,9).Value = 1594.7059
$ws.Cells.Item(99,10).Value = 2552.75
$ws.Cells.Item(99,11).Value = 1594.7059
$ws.Cells.Item(99,12).Value = 2552.75
$ws.Cells.Item(99,13).Value = -96.70589999999993
$ws.Cells.Item(99,14).Value = -5548.75

$ws.Cells.Item(105,8).Value = 5010
$ws.Cells.Item(105,9).Value = 5010
$ws.Cells.Item(105,10).Value = 0
$ws.Cells.Item(105,11).Value = 5010
$ws.Cells.Item(105,12).Value = 0
$ws.Cells.Item(105,13).Value = -3263
$ws.Cells.Item(105,14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99,8).Value = 3575.3333
$ws.Cells.Item(99,9).Value = 3604.8096
$ws.Cells.Item(99,10).Value = 3472.1667
$ws.Cells.Item(99,11).Value = 3604.8096
$ws.Cells.Item(99,12).Value = 3472.1667
$ws.Cells.Item(99,13).Value = -2106.8096
$ws.Cells.Item(99,14).Value = -6468.1667

$ws.Cells.Item(107,8).Value = 404.07407
$ws.Cells.Item(107,9).Value = 351
$ws.Cells.Item(107,10).Value = 481.27274
$ws.Cells.Item(107,11).Value = 351
$ws.Cells.Item(107,12).Value = 481.27274
$ws.Cells.Item(107,13).Value = 1569
$ws.Cells.Item(107,14).Value = -4321.27274

$ws.Cells.Item(126,8).Value = 3575.3333
$ws.Cells.Item(126,9).Value = 3604.8096
$ws.Cells.Item(126,10).Value = 3472.1667
$ws.Cells.Item(126,11).Value = 10814.4288
$ws.Cells.Item(126,12).Value = 10416.5001
$ws.Cells.Item(126,13).Value = -8344.4288
$ws.Cells.Item(126,14).Value = -15356.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23,8).Value = 380
$ws.Cells.Item(23,10).Value = 340
$ws.Cells.Item(23,12).Value = 1020
$ws.Cells.Item(23,14).Value = -1490

$ws.Cells.Item(75,8).Value = 5055.727
$ws.Cells.Item(75,9).Value = 2478.25
$ws.Cells.Item(75,10).Value = 6528.5713
$ws.Cells.Item(75,11).Value = 7434.75
$ws.Cells.Item(75,12).Value = 19585.7139
$ws.Cells.Item(75,13).Value = -6436.75
$ws.Cells.Item(75,14).Value = -21581.7139

$ws.Cells.Item(78,8).Value = 5055.727
$ws.Cells.Item(78,9).Value = 2478.25
$ws.Cells.Item(78,10).Value = 6528.5713
$ws.Cells.Item(78,11).Value = 22304.25
$ws.Cells.Item(78,12).Value = 58757.14169999999
$ws.Cells.Item(78,13).Value = -17312.25
$ws.Cells.Item(78,14).Value = -68741.14169999999

$ws.Cells.Item(113,8).Value = 680.94446
$ws.Cells.Item(113,9).Value = 679.2
$ws.Cells.Item(113,10).Value = 681.61536
$ws.Cells.Item(113,11).Value = 2037.6
$ws.Cells.Item(113,12).Value = 2044.84608
$ws.Cells.Item(113,13).Value = 132.3999999999999
$ws.Cells.Item(113,14).Value = -6384.84608

$ws.Cells.Item(133,8).Value = 4100.909
$ws.Cells.Item(133,9).Value = 1858
$ws.Cells.Item(133,10).Value = 5970
$ws.Cells.Item(133,11).Value = 5574
$ws.Cells.Item(133,12).Value = 17910
$ws.Cells.Item(133,13).Value = -514
$ws.Cells.Item(133,14).Value = -28030

$ws.Cells.Item(140,8).Value = 2556.647
$ws.Cells.Item(140,9).Value = 1010
$ws.Cells.Item(140,10).Value = 3931.4443
$ws.Cells.Item(140,11).Value = 3030
$ws.Cells.Item(140,12).Value = 11794.3329
$ws.Cells.Item(140,13).Value = 2150
$ws.Cells.Item(140,14).Value = -22154.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22,8).Value = 622.2222
$ws.Cells.Item(22,9).Value = 662.5
$ws.Cells.Item(22,10).Value = 300
$ws.Cells.Item(22,11).Value = 662.5
$ws.Cells.Item(22,12).Value = 300
$ws.Cells.Item(22,13).Value = -367.5
$ws.Cells.Item(22,14).Value = -890

$ws.Cells.Item(27,8).Value = 622.2222
$ws.Cells.Item(27,9).Value = 662.5
$ws.Cells.Item(27,10).Value = 300
$ws.Cells.Item(27,11).Value = 662.5
$ws.Cells.Item(27,12).Value = 300
$ws.Cells.Item(27,13).Value = -555.5
$ws.Cells.Item(27,14).Value = -514

$ws.Cells.Item(76,8).Value = 24000
$ws.Cells.Item(76,10).Value = 24000
$ws.Cells.Item(76,12).Value = 24000
$ws.Cells.Item(76,14).Value = -24676

$ws.Cells.Item(79,8).Value = 24000
$ws.Cells.Item(79,10).Value = 24000
$ws.Cells.Item(79,12).Value = 24000
$ws.Cells.Item(79,14).Value = -26340

$ws.Cells.Item(95,8).Value = 11500
$ws.Cells.Item(95,10).Value = 11500
$ws.Cells.Item(95,12).Value = 11500
$ws.Cells.Item(95,14).Value = -16992

$ws.Cells.Item(100,8).Value = 7466.6665
$ws.Cells.Item(100,9).Value = 17666.666
$ws.Cells.Item(100,10).Value = 2366.6667
$ws.Cells.Item(100,11).Value = 17666.666
$ws.Cells.Item(100,12).Value = 2366.6667
$ws.Cells.Item(100,13).Value = -17125.666
$ws.Cells.Item(100,14).Value = -3448.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54,8).Value = 26075
$ws.Cells.Item(54,10).Value = 26075
$ws.Cells.Item(54,12).Value = 26075
$ws.Cells.Item(54,14).Value = -27115

$ws.Cells.Item(97,8).Value = 23286
$ws.Cells.Item(97,10).Value = 23286
$ws.Cells.Item(97,12).Value = 23286
$ws.Cells.Item(97,14).Value = -25268

$ws.Cells.Item(107,8).Value = 564.4
$ws.Cells.Item(107,9).Value = 505.5
$ws.Cells.Item(107,10).Value = 800
$ws.Cells.Item(107,11).Value = 1516.5
$ws.Cells.Item(107,12).Value = 2400
$ws.Cells.Item(107,13).Value = 403.5
$ws.Cells.Item(107,14).Value = -6240

$ws.Cells.Item(113,8).Value = 1058
$ws.Cells.Item(113,9).Value = 841
$ws.Cells.Item(113,10).Value = 1323.2222
$ws.Cells.Item(113,11).Value = 2523
$ws.Cells.Item(113,12).Value = 3969.6666
$ws.Cells.Item(113,13).Value = -353
$ws.Cells.Item(113,14).Value = -8309.6666

$ws.Cells.Item(123,8).Value = 57193.707
$ws.Cells.Item(123,10).Value = 57193.707
$ws.Cells.Item(123,12).Value = 57193.707
$ws.Cells.Item(123,14).Value = -66993.70699999999
